$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 6.15379541431027),
    @(3, 0.6545652718822623, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 2.964545797025059),
    @(4, 0.6545652718822623, 0.3048912486333797, 0.7210945179870265, 13.86384647080068, 1, 15.54439750930335),
    @(5, 0.6545652718822623, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 2.964545797025059),
    @(6, 0.6545652718822623, 1.626987699542094, 0.7210945179870265, 13.86384647080068, 1, 16.86649396021207),
    @(7, 3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 1, 8.656069925401464),
    @(8, 1.445647641019636, 0.04103571897497393, 0.7210945179870265, 0.5333859586016987, 1, 2.741163836583335),
    @(9, 3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 5.582307763322248),
    @(10, 3.272327238179451, 1.626987699542094, 3.223369029078222, 13.86384647080068, 0, 21.98653043760045),
    @(11, 1.445647641019636, 1.626987699542094, 0.7210945179870265, 13.86384647080068, 1, 17.65757632934944),
    @(12, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 6.15379541431027)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
